$d = $word.ActiveDocument

# 1) Update the heading text "Requisitos Funcionais:" -> "Requisitos do Sistema:"
$d.Content.Find.Execute("Requisitos Funcionais:", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Requisitos do Sistema:", 2)

# 2) Resize the table's two grid columns: 1230 -> 1185 and 7785 -> 7830 (twips)
#    Column.Width is expressed in points (1 pt = 20 twips) in the Word OM.
$table = $d.Tables(1)
$table.Columns(1).Width = 1185 / 20
$table.Columns(2).Width = 7830 / 20
